$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-5, per diff
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 5
